# Applies the commit "Add data for 2021-12-02":
#  - Rename the sheet / update the "through" date from Nov 23 to Nov 24
#  - Update the header shared string text to match
#  - Bump a handful of monthly neighborhood carjacking counts (and add
#    a few brand-new counts) reflecting one additional day of data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename worksheet tab (Through 2021-11-23 -> Through 2021-11-24)
$ws.Name = "Through 2021-11-24"

# Update the column header text for the "current month" column
$ws.Range("B1").Value = "November 2021 (through November 24)"

# Cell value updates (existing cells incremented by 1)
$ws.Range("M2").Value = 16
$ws.Range("M3").Value = 15
$ws.Range("M4").Value = 11
$ws.Range("AT4").Value = 11
$ws.Range("BE5").Value = 6
$ws.Range("BE7").Value = 3
$ws.Range("B8").Value = 6
$ws.Range("X8").Value = 3
$ws.Range("AI8").Value = 3
$ws.Range("M11").Value = 3
$ws.Range("M13").Value = 4
$ws.Range("B19").Value = 2
$ws.Range("BE20").Value = 2
$ws.Range("X24").Value = 3
$ws.Range("M28").Value = 2
$ws.Range("M31").Value = 5
$ws.Range("AT42").Value = 2
$ws.Range("M49").Value = 2
$ws.Range("M55").Value = 2
$ws.Range("B68").Value = 4
$ws.Range("B80").Value = 4

# Newly added cell values (previously blank)
$ws.Range("M14").Value = 1
$ws.Range("B32").Value = 1
$ws.Range("X34").Value = 1
$ws.Range("BP48").Value = 1
$ws.Range("M52").Value = 1
$ws.Range("B67").Value = 1
$ws.Range("AI72").Value = 1
$ws.Range("M73").Value = 1
